$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for each coin row.
# D-column values are entered with a leading apostrophe so Excel keeps
# them as literal text (matching the original inlineStr cells) instead of
# re-interpreting numeric-looking strings as numbers (which would drop
# formatting like trailing zeros or thousands-separator dots). The
# Style reset afterwards clears the "quote prefix" cell style Excel
# applies automatically, so the cell keeps its original (default) style.
$ws.Range("D2").Value = "'92.001.76"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.21%  "
$ws.Range("D3").Value = "'3.327.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.68%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'230.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.66%  "
$ws.Range("D6").Value = "'613.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.63%  "
$ws.Range("E7").Value = "  -1.18%  "
$ws.Range("D8").Value = "'0.384"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.29%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").Value = "'0.958"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.26%  "
$ws.Range("D11").Value = "'3.327.15"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.68%  "
$ws.Range("D12").Value = "'42.60"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.83%  "
$ws.Range("E13").Value = "  -1.17%  "
$ws.Range("E14").Value = "  +0.89%  "
$ws.Range("D15").Value = "'91.815.29"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.06%  "
$ws.Range("D16").Value = "'3.951.24"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.58%  "
$ws.Range("E17").Value = "  -1.93%  "
$ws.Range("D18").Value = "'8.04"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.45%  "
$ws.Range("D19").Value = "'3.326.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.73%  "
$ws.Range("D20").Value = "'17.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.69%  "
$ws.Range("D21").Value = "'10.83"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.76%  "
$ws.Range("D22").Value = "'3.43"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.56%  "
$ws.Range("D23").Value = "'491.44"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").Value = "'0.437"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -9.89%  "
$ws.Range("D25").Value = "'6.55"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.29%  "
$ws.Range("E26").Value = "  -4.29%  "
$ws.Range("D27").Value = "'92.55"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.41%  "
$ws.Range("D28").Value = "'11.90"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.65%  "
$ws.Range("D29").Value = "'3.505.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.64%  "
$ws.Range("E30").Value = "  +0.45%  "
$ws.Range("D31").Value = "'11.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.53%  "
$ws.Range("E32").Value = "  +2.59%  "
$ws.Range("E33").Value = "  -3.78%  "
$ws.Range("E34").Value = "  +0.17%  "
$ws.Range("D35").Value = "'0.174"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.95%  "
$ws.Range("D36").Value = "'28.26"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.16%  "
$ws.Range("E37").Value = "  -3.72%  "
$ws.Range("D38").Value = "'560.74"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.75%  "
$ws.Range("D39").Value = "'7.42"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.17%  "
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("E41").Value = "  -0.48%  "
$ws.Range("D42").Value = "'1.37"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.16%  "
$ws.Range("E43").Value = "  -4.08%  "
$ws.Range("D44").Value = "'23.71"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.06%  "
$ws.Range("D45").Value = "'0.0414"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.28%  "
$ws.Range("D46").Value = "'1.67"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.19%  "
$ws.Range("D47").Value = "'3.60"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.28%  "
$ws.Range("E48").Value = "  -0.91%  "
$ws.Range("E49").Value = "  -0.33%  "
$ws.Range("D50").Value = "'8.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.76%  "
$ws.Range("D51").Value = "'52.05"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.26%  "
